# 2012 Volunteer Coordination Plan - apply commit "Updated the mailer and doc"
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Collapse(1)
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: could not find text: $findText"
    } else {
        Write-Host "OK: replaced '$findText' -> '$replaceText'"
    }
}

# 1) "...raising 45,000 pounds of food" -> "...raising 25,000 pounds of food"
Replace-Text "45,000 pounds of food" "25,000 pounds of food"

# 2) Week prior to the event date range: "2-9, 2011" -> "3-7, 2012", plus a _GoBack bookmark
#    right after the new date (mirrors Word's automatic "last edit" bookmark).
Replace-Text "2-9, 2011" "3-7, 2012"

$rngAfterWeek = $d.Content
$rngAfterWeek.Collapse(1)
if ($rngAfterWeek.Find.Execute("3-7, 2012")) {
    $goBackRng = $d.Range($rngAfterWeek.End, $rngAfterWeek.End)
    $d.Bookmarks.Add("_GoBack", $goBackRng)
    Write-Host "OK: added _GoBack bookmark"
} else {
    Write-Host "WARNING: could not locate new week-prior date to place bookmark"
}

# 3) Bag Pickup paragraph: add the RightNow / Prudential pickup locations after "...for pickup"
$rngPickup = $d.Content
$rngPickup.Collapse(1)
if ($rngPickup.Find.Execute("Gallatin Valley Food bank for pickup")) {
    $insertPt = $d.Range($rngPickup.End, $rngPickup.End)
    $insertPt.InsertAfter(", the former RightNow Technologies office locations and the Prudential office on Stadium Drive")
    Write-Host "OK: inserted pickup locations sentence"
} else {
    Write-Host "WARNING: could not find Gallatin Valley Food bank for pickup"
}

# 4) Day of event date: "10, 2011" -> "8, 2012"
Replace-Text " 10, 2011" " 8, 2012"

# 5) Remove the "Marty: 406-570-0111" list paragraph entirely (Ron's paragraph moves up)
$deleted = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Marty*406-570-0111*") {
        $para.Range.Delete()
        $deleted = $true
        break
    }
}
if ($deleted) {
    Write-Host "OK: removed Marty paragraph"
} else {
    Write-Host "WARNING: Marty paragraph not found"
}

# 6) Drop-off location list: "TBD – Probably Rosauers or the Main Mall" -> "Rosauers "
$rngTbd = $d.Content
$rngTbd.Collapse(1)
if ($rngTbd.Find.Execute("TBD – Probably ")) {
    $rngTbd.Text = ""
    Write-Host "OK: removed 'TBD - Probably ' prefix"
} else {
    Write-Host "WARNING: 'TBD - Probably ' prefix not found"
}

$rngMall = $d.Content
$rngMall.Collapse(1)
if ($rngMall.Find.Execute(" or the Main Mall")) {
    $rngMall.Text = " "
    Write-Host "OK: replaced ' or the Main Mall' with a single space"
} else {
    Write-Host "WARNING: ' or the Main Mall' not found"
}

# 7) Update the closing mailto hyperlink from Marty to Kristen
if ($d.Hyperlinks.Count -ge 1) {
    $h = $d.Hyperlinks(1)
    $h.TextToDisplay = "Kristen.radford@oracle.com"
    $h.Address = "mailto:Kristen.radford@oracle.com"
    Write-Host "OK: updated hyperlink to Kristen.radford@oracle.com"
} else {
    Write-Host "WARNING: no hyperlinks found in document"
}
